$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the "Current Status" for row 10 (ASP .Net MVC) to "In progress"
$ws.Range("C10").Value = "In progress"

# Move the active selection to C11, matching the workbook's saved UI state
$ws.Range("C11").Select()
